$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 300547.53
$ws.Range("I33").Value = 382.6842
$ws.Range("J33").Value = 1115280.8
$ws.Range("K33").Value = 382.6842
$ws.Range("L33").Value = 1115280.8
$ws.Range("M33").Value = -153.6842
$ws.Range("N33").Value = -1115738.8
$ws.Range("H76").Value = 4459.593
$ws.Range("I76").Value = 4364.0454
$ws.Range("K76").Value = 4364.0454
$ws.Range("M76").Value = -4049.0454
$ws.Range("H79").Value = 4459.593
$ws.Range("I79").Value = 4364.0454
$ws.Range("K79").Value = 4364.0454
$ws.Range("M79").Value = -3272.0454
$ws.Range("H86").Value = 38922.223
$ws.Range("J86").Value = 1835.8889
$ws.Range("L86").Value = 1835.8889
$ws.Range("N86").Value = -4081.8889
$ws.Range("H89").Value = 38922.223
$ws.Range("J89").Value = 1835.8889
$ws.Range("L89").Value = 9179.4445
$ws.Range("N89").Value = -20411.4445
$ws.Range("H112").Value = 1019.24
$ws.Range("I112").Value = 490
$ws.Range("J112").Value = 1041.2916
$ws.Range("K112").Value = 1470
$ws.Range("L112").Value = 3123.8748
$ws.Range("M112").Value = -362
$ws.Range("N112").Value = -5339.8748
$ws.Range("H138").Value = 2795.0942
$ws.Range("I138").Value = 1870.25
$ws.Range("J138").Value = 3195.027
$ws.Range("K138").Value = 5610.75
$ws.Range("L138").Value = 9585.081
$ws.Range("M138").Value = -470.75
$ws.Range("N138").Value = -19865.081

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10163.293
$ws.Range("I32").Value = 8918.057000000001
$ws.Range("K32").Value = 8918.057000000001
$ws.Range("M32").Value = -8631.057000000001
$ws.Range("H74").Value = 1596.421
$ws.Range("I74").Value = 1873.5
$ws.Range("K74").Value = 1873.5
$ws.Range("M74").Value = -999.5
$ws.Range("H77").Value = 1596.421
$ws.Range("I77").Value = 1873.5
$ws.Range("K77").Value = 9367.5
$ws.Range("M77").Value = -4999.5
$ws.Range("H132").Value = 4936.909
$ws.Range("I132").Value = 5004.44
$ws.Range("J132").Value = 4725.875
$ws.Range("K132").Value = 15013.32
$ws.Range("L132").Value = 14177.625
$ws.Range("M132").Value = -12483.32
$ws.Range("N132").Value = -19237.625

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2800.7932
$ws.Range("I134").Value = 2598.818
$ws.Range("J134").Value = 3435.5715
$ws.Range("K134").Value = 7796.454000000001
$ws.Range("L134").Value = 10306.7145
$ws.Range("M134").Value = -5261.454000000001
$ws.Range("N134").Value = -15376.7145

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1734.6666
$ws.Range("I3").Value = 102
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 102
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = 11
$ws.Range("N3").Value = -5226
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H134").Value = 2416.6667
$ws.Range("I134").Value = 1500
$ws.Range("K134").Value = 4500
$ws.Range("M134").Value = -1965

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 101.666664
$ws.Range("I10").Value = 101.666664
$ws.Range("K10").Value = 304.999992
$ws.Range("M10").Value = -165.999992
$ws.Range("H28").Value = 1147
$ws.Range("I28").Value = 1007.5
$ws.Range("J28").Value = 1258.6
$ws.Range("K28").Value = 3022.5
$ws.Range("L28").Value = 3775.8
$ws.Range("M28").Value = -2790.5
$ws.Range("N28").Value = -4239.799999999999
$ws.Range("H39").Value = 500
$ws.Range("I39").Value = 500
$ws.Range("K39").Value = 1500
$ws.Range("M39").Value = -1206
$ws.Range("H45").Value = 1264.5
$ws.Range("I45").Value = 800
$ws.Range("J45").Value = 1357.4
$ws.Range("K45").Value = 2400
$ws.Range("L45").Value = 4072.2
$ws.Range("M45").Value = -1868
$ws.Range("N45").Value = -5136.200000000001
$ws.Range("H52").Value = 1895.3334
$ws.Range("J52").Value = 1895.3334
$ws.Range("L52").Value = 5686.0002
$ws.Range("N52").Value = -6218.0002
$ws.Range("H112").Value = 2320317.5
$ws.Range("I112").Value = 333999.34
$ws.Range("K112").Value = 1001998.02
$ws.Range("M112").Value = -1000890.02
$ws.Range("H113").Value = 729.8222
$ws.Range("I113").Value = 1117.4375
$ws.Range("J113").Value = 515.9655
$ws.Range("K113").Value = 3352.3125
$ws.Range("L113").Value = 1547.8965
$ws.Range("M113").Value = -1182.3125
$ws.Range("N113").Value = -5887.8965
$ws.Range("H114").Value = 1229.3
$ws.Range("I114").Value = 540
$ws.Range("J114").Value = 1688.8334
$ws.Range("K114").Value = 1620
$ws.Range("L114").Value = 5066.5002
$ws.Range("M114").Value = 1634
$ws.Range("N114").Value = -11574.5002
$ws.Range("H116").Value = 2419.6
$ws.Range("I116").Value = 3499
$ws.Range("J116").Value = 1700
$ws.Range("K116").Value = 10497
$ws.Range("L116").Value = 5100
$ws.Range("M116").Value = -7055
$ws.Range("N116").Value = -11984
$ws.Range("H118").Value = 2441.3333
$ws.Range("I118").Value = 882.6667
$ws.Range("J118").Value = 4000
$ws.Range("K118").Value = 2648.0001
$ws.Range("L118").Value = 12000
$ws.Range("M118").Value = -1405.0001
$ws.Range("N118").Value = -14486
$ws.Range("H119").Value = 333750
$ws.Range("I119").Value = 333750
$ws.Range("K119").Value = 1001250
$ws.Range("M119").Value = -996412
$ws.Range("H120").Value = 745209
$ws.Range("I120").Value = 745209
$ws.Range("K120").Value = 2235627
$ws.Range("M120").Value = -2230789
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H126").Value = 2404.2856
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -18880
$ws.Range("H129").Value = 215060.12
$ws.Range("I129").Value = 7605.75
$ws.Range("J129").Value = 318787.3
$ws.Range("K129").Value = 22817.25
$ws.Range("L129").Value = 956361.8999999999
$ws.Range("M129").Value = -17817.25
$ws.Range("N129").Value = -966361.8999999999
$ws.Range("H131").Value = 822.79
$ws.Range("J131").Value = 827.33673
$ws.Range("L131").Value = 2482.01019
$ws.Range("N131").Value = -12562.01019
$ws.Range("H136").Value = 1757.375
$ws.Range("I136").Value = 1346.5
$ws.Range("K136").Value = 4039.5
$ws.Range("M136").Value = 1060.5
$ws.Range("H138").Value = 1344.8
$ws.Range("I138").Value = 1344.8
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 4034.4
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 1105.6
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 1985.6757
$ws.Range("I139").Value = 742.5
$ws.Range("J139").Value = 2582.4
$ws.Range("K139").Value = 2227.5
$ws.Range("L139").Value = 7747.200000000001
$ws.Range("M139").Value = 2912.5
$ws.Range("N139").Value = -18027.2
$ws.Range("H140").Value = 6729.2104
$ws.Range("I140").Value = 6986.3887
$ws.Range("K140").Value = 20959.1661
$ws.Range("M140").Value = -15779.1661

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 42132.832
$ws.Range("J69").Value = 42132.832
$ws.Range("L69").Value = 42132.832
$ws.Range("N69").Value = -43630.832
$ws.Range("H72").Value = 42132.832
$ws.Range("J72").Value = 42132.832
$ws.Range("L72").Value = 126398.496
$ws.Range("N72").Value = -133886.496

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1762
$ws.Range("I93").Value = 1715.2
$ws.Range("J93").Value = 1847.091
$ws.Range("K93").Value = 1715.2
$ws.Range("L93").Value = 1847.091
$ws.Range("M93").Value = -467.2
$ws.Range("N93").Value = -4343.091

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1942.7301
$ws.Range("I132").Value = 1793.6078
$ws.Range("J132").Value = 2576.5
$ws.Range("K132").Value = 5380.8234
$ws.Range("L132").Value = 7729.5
$ws.Range("M132").Value = -2850.8234
$ws.Range("N132").Value = -12789.5
